$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(27, 9).Value = 'sd'
$ws.Cells.Item(27, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(40, 9).Value = 'sd'
$ws.Cells.Item(40, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(43, 9).Value = 'sd'
$ws.Cells.Item(43, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(52, 9).Value = 'sd'
$ws.Cells.Item(52, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(77, 9).Value = '%'
$ws.Cells.Item(77, 10).Value = 'Uninterpretable'
$ws.Cells.Item(92, 9).Value = 'sd'
$ws.Cells.Item(92, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(94, 9).Value = 'sd'
$ws.Cells.Item(94, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(108, 9).Value = 'sd'
$ws.Cells.Item(108, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(110, 9).Value = 'sd'
$ws.Cells.Item(110, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(113, 9).Value = 'sd'
$ws.Cells.Item(113, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(123, 9).Value = 'aa'
$ws.Cells.Item(123, 10).Value = 'Agree/Accept'
$ws.Cells.Item(124, 9).Value = 'aa'
$ws.Cells.Item(124, 10).Value = 'Agree/Accept'
$ws.Cells.Item(131, 9).Value = 'sv'
$ws.Cells.Item(131, 10).Value = 'Statement-opinion'
$ws.Cells.Item(136, 9).Value = 'sv'
$ws.Cells.Item(136, 10).Value = 'Statement-opinion'
$ws.Cells.Item(146, 9).Value = 'aa'
$ws.Cells.Item(146, 10).Value = 'Agree/Accept'
$ws.Cells.Item(166, 9).Value = 'qy'
$ws.Cells.Item(166, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(177, 9).Value = 'b'
$ws.Cells.Item(177, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(184, 9).Value = 'sv'
$ws.Cells.Item(184, 10).Value = 'Statement-opinion'
$ws.Cells.Item(191, 9).Value = 'ba'
$ws.Cells.Item(191, 10).Value = 'Appreciation'
$ws.Cells.Item(204, 9).Value = 'sd'
$ws.Cells.Item(204, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(209, 9).Value = 'sd'
$ws.Cells.Item(209, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(214, 9).Value = 'sd'
$ws.Cells.Item(214, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(222, 9).Value = 'aa'
$ws.Cells.Item(222, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'aa'
$ws.Cells.Item(223, 10).Value = 'Agree/Accept'
$ws.Cells.Item(225, 9).Value = 'sd'
$ws.Cells.Item(225, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(229, 9).Value = 'b'
$ws.Cells.Item(229, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(232, 9).Value = 'b'
$ws.Cells.Item(232, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(236, 9).Value = 'sd'
$ws.Cells.Item(236, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(242, 9).Value = 'sv'
$ws.Cells.Item(242, 10).Value = 'Statement-opinion'
$ws.Cells.Item(249, 9).Value = 'aa'
$ws.Cells.Item(249, 10).Value = 'Agree/Accept'
$ws.Cells.Item(250, 9).Value = 'sd'
$ws.Cells.Item(250, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(252, 9).Value = 'aa'
$ws.Cells.Item(252, 10).Value = 'Agree/Accept'
$ws.Cells.Item(259, 9).Value = 'aa'
$ws.Cells.Item(259, 10).Value = 'Agree/Accept'
$ws.Cells.Item(262, 9).Value = 'sd'
$ws.Cells.Item(262, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(302, 9).Value = 'sd'
$ws.Cells.Item(302, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(304, 9).Value = 'sd'
$ws.Cells.Item(304, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(306, 9).Value = 'aa'
$ws.Cells.Item(306, 10).Value = 'Agree/Accept'
$ws.Cells.Item(314, 9).Value = 'sd'
$ws.Cells.Item(314, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(315, 9).Value = 'b'
$ws.Cells.Item(315, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(318, 9).Value = '%'
$ws.Cells.Item(318, 10).Value = 'Uninterpretable'
$ws.Cells.Item(326, 9).Value = '%'
$ws.Cells.Item(326, 10).Value = 'Uninterpretable'
$ws.Cells.Item(331, 9).Value = '%'
$ws.Cells.Item(331, 10).Value = 'Uninterpretable'
$ws.Cells.Item(343, 9).Value = 'ba'
$ws.Cells.Item(343, 10).Value = 'Appreciation'
$ws.Cells.Item(344, 9).Value = 'sv'
$ws.Cells.Item(344, 10).Value = 'Statement-opinion'
$ws.Cells.Item(348, 9).Value = 'b'
$ws.Cells.Item(348, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(353, 9).Value = 'b'
$ws.Cells.Item(353, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(355, 9).Value = 'b'
$ws.Cells.Item(355, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(359, 9).Value = 'b'
$ws.Cells.Item(359, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(362, 9).Value = 'sd'
$ws.Cells.Item(362, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(363, 9).Value = 'qy'
$ws.Cells.Item(363, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(378, 9).Value = 'b'
$ws.Cells.Item(378, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(388, 9).Value = 'ba'
$ws.Cells.Item(388, 10).Value = 'Appreciation'
$ws.Cells.Item(392, 9).Value = 'sd'
$ws.Cells.Item(392, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(394, 9).Value = 'sv'
$ws.Cells.Item(394, 10).Value = 'Statement-opinion'
$ws.Cells.Item(406, 9).Value = 'sd'
$ws.Cells.Item(406, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(438, 9).Value = 'b'
$ws.Cells.Item(438, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(454, 9).Value = 'sd'
$ws.Cells.Item(454, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(468, 9).Value = 'ba'
$ws.Cells.Item(468, 10).Value = 'Appreciation'
$ws.Cells.Item(470, 9).Value = 'sd'
$ws.Cells.Item(470, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(476, 9).Value = 'ba'
$ws.Cells.Item(476, 10).Value = 'Appreciation'
$ws.Cells.Item(504, 9).Value = 'sv'
$ws.Cells.Item(504, 10).Value = 'Statement-opinion'
$ws.Cells.Item(509, 9).Value = 'sd'
$ws.Cells.Item(509, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(511, 9).Value = 'sv'
$ws.Cells.Item(511, 10).Value = 'Statement-opinion'
$ws.Cells.Item(518, 9).Value = 'sd'
$ws.Cells.Item(518, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(527, 9).Value = 'sd'
$ws.Cells.Item(527, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(543, 9).Value = 'sv'
$ws.Cells.Item(543, 10).Value = 'Statement-opinion'
$ws.Cells.Item(546, 9).Value = 'sv'
$ws.Cells.Item(546, 10).Value = 'Statement-opinion'
$ws.Cells.Item(552, 9).Value = 'qy'
$ws.Cells.Item(552, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(553, 9).Value = '%'
$ws.Cells.Item(553, 10).Value = 'Uninterpretable'
$ws.Cells.Item(555, 9).Value = 'sv'
$ws.Cells.Item(555, 10).Value = 'Statement-opinion'
$ws.Cells.Item(557, 9).Value = 'sd'
$ws.Cells.Item(557, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(558, 9).Value = 'sd'
$ws.Cells.Item(558, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(559, 9).Value = 'sv'
$ws.Cells.Item(559, 10).Value = 'Statement-opinion'
$ws.Cells.Item(562, 9).Value = 'aa'
$ws.Cells.Item(562, 10).Value = 'Agree/Accept'
$ws.Cells.Item(565, 9).Value = 'sd'
$ws.Cells.Item(565, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(585, 9).Value = 'sd'
$ws.Cells.Item(585, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(592, 9).Value = 'fc'
$ws.Cells.Item(592, 10).Value = 'Conventional-closing'
